# Insert a new slide ("Visualisation de données") as the new slide 6,
# pushing the former slides 6 ("modélisation") and 7 ("Gains de
# performance des différents modèles") down to 7 and 8 respectively.

$p = $ppt.ActivePresentation

# Slide 5 ("Récupération de données") uses the "Titre et contenu" layout
# (slideLayout2.xml) -- legacy layout index 2 maps to the same layout, so
# the new slide matches the look of its neighbours.
$newSlide = $p.Slides.Add(6, 2)

# Turn on the slide-number placeholder (the new slide mirrors slide 5,
# which shows the slide number in the bottom-right corner).
$newSlide.HeadersFooters.SlideNumber.Visible = $true

# Match the naming convention ("Titre N" / "Espace réservé du contenu N" /
# "Espace réservé du numéro de diapositive N") used by the rest of the deck.
$newSlide.Shapes.Item(1).Name = "Titre 1"
$newSlide.Shapes.Item(2).Name = "Espace réservé du contenu 2"
$newSlide.Shapes.Item(3).Name = "Espace réservé du numéro de diapositive 4"

# Title placeholder.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Visualisation de données"

# Body / content placeholder.
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Nous avons affiché quelques histogrammes relatifs à l’UPDRS "
